# Add PF/1.0.3 to meta-sheet
# Appends a new row (row 3) below the existing header (row 1) and
# PF/1.0.0 row (row 2), recording the new release "PF/1.0.3" together
# with "X" markers in the remaining columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.3"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# The new row carries no special formatting (unlike rows 1-2, which use
# the sheet's bold/sized header style), so reset it back to the default
# "Normal" style.
$ws.Range("A3:D3").Style = "Normal"
